# The workbook's "Oman MSME" summary sheet lists, under the MSMEs column,
# two stat rows back to back:
#   row 12: Enterprises (absolute #)            118386
#   row 13: Enterprises density (per 1000 people) 44.5
#
# The commit reorders these two rows so the "density" metric is reported
# before the "absolute #" metric:
#   row 12: Enterprises density (per 1000 people) 44.5
#   row 13: Enterprises (absolute #)             118386
#
# Swap the label/value pairs between the two rows. The D-column values are
# numeric-looking text (stored as strings, not numbers) in the source file,
# so force a text number format before assigning them to keep them as text
# instead of being auto-coerced into numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Enterprises density (per 1000 people)"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.5"

$ws.Range("A13").Value = "Enterprises (absolute #)"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "118386"
